# Update "想去人数" (want-to-go count) values across the four sheets to
# reflect a refreshed data export (gh-pages output regenerated).
#
# Sheet order in the workbook:
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life)
#   4 = 全部类型  (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F7").Value  = 5575
$ws1.Range("F10").Value = 3817
$ws1.Range("F20").Value = 463
$ws1.Range("F21").Value = 128
$ws1.Range("F23").Value = 5241
$ws1.Range("F24").Value = 437
$ws1.Range("F25").Value = 2079
$ws1.Range("F28").Value = 7774
$ws1.Range("F31").Value = 2188
$ws1.Range("F32").Value = 2165
$ws1.Range("F33").Value = 1328
$ws1.Range("F37").Value = 18
$ws1.Range("F38").Value = 264
$ws1.Range("F43").Value = 1172
$ws1.Range("F45").Value = 1323
$ws1.Range("F46").Value = 2048
$ws1.Range("F48").Value = 215
$ws1.Range("F49").Value = 1213

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 119

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 731

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value  = 731
$ws4.Range("F9").Value  = 5575
$ws4.Range("F10").Value = 3817
$ws4.Range("F20").Value = 463
$ws4.Range("F22").Value = 128
$ws4.Range("F24").Value = 5241
$ws4.Range("F25").Value = 437
$ws4.Range("F26").Value = 2079
$ws4.Range("F29").Value = 7774
$ws4.Range("F32").Value = 2188
$ws4.Range("F33").Value = 2165
$ws4.Range("F34").Value = 1328
$ws4.Range("F37").Value = 264
$ws4.Range("F41").Value = 1172
$ws4.Range("F43").Value = 1323
$ws4.Range("F45").Value = 2048
$ws4.Range("F48").Value = 215
$ws4.Range("F49").Value = 1213
